# Att - Antes do almoco
# Append new registry rows (VDS1842..VDS1850) to the "Sheet" worksheet,
# right after the existing data that ends at row 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("VDS1842", "15/10/2025 11:01:13", "1585"),
    @("VDS1843", "15/10/2025 11:05:59", "555555555"),
    @("VDS1844", "15/10/2025 11:06:15", "555555555"),
    @("VDS1845", "15/10/2025 11:06:34", "555555555"),
    @("VDS1846", "15/10/2025 11:06:56", "555555555"),
    @("VDS1847", "15/10/2025 11:08:04", "555555555"),
    @("VDS1848", "15/10/2025 11:08:32", "555555555"),
    @("VDS1849", "15/10/2025 11:34:44", "123563"),
    @("VDS1850", "15/10/2025 11:34:50", "123563")
)

$startRow = 38
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $codigo = $rows[$i][0]
    $dataRegistro = $rows[$i][1]
    $projeto = $rows[$i][2]

    # Column A/B are plain, non-numeric text - assigning the string is enough
    # to keep them stored as text.
    $ws.Cells.Item($r, 1).Value = $codigo
    $ws.Cells.Item($r, 2).Value = $dataRegistro

    # Column C values are digit-only strings (e.g. "1585", "555555555").
    # Excel auto-converts a plain digit string to a number, so force text
    # with a leading apostrophe, then reset the cell style back to Normal
    # so no stray "quote prefix" number format sticks around.
    $ws.Cells.Item($r, 3).Value = "'" + $projeto
    $ws.Cells.Item($r, 3).Style = "Normal"
}
